# estructurasDatosWordix.xlsx edit
# Adds new words to $coleccionPalabras and documents a new associative
# array $coleccionPartidas, per commit:
#   "Agregue palabras al arreglo de palabras, modifique estructurasDatosWordix"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Clear the cells from the old layout that are being replaced/relocated
# ---------------------------------------------------------------------
$ws.Range("A3:H4").ClearContents()
$ws.Range("B6:B9").ClearContents()
$ws.Cells.Item(13, 1).ClearContents()

# ---------------------------------------------------------------------
# Title
# ---------------------------------------------------------------------
$ws.Cells.Item(1, 1).Value = "GRUPO 8: Cruz Oviedo, De Phillipis, Padilla"

# ---------------------------------------------------------------------
# $coleccionPalabras section
# ---------------------------------------------------------------------
$ws.Cells.Item(3, 1).Value = "Representacion del arreglo `$coleccionPalabras"
# row 4 left blank (styled header continuation)

$ws.Cells.Item(5, 1).Value = "`$coleccionPalabras="

$indices = 0,1,2,3,4,5,6,7,8,9,10
$col = 4
foreach ($i in $indices) {
    $ws.Cells.Item(5, $col).Value = $i
    $col++
}
$ws.Cells.Item(5, 15).Value = "//indices"

$palabras = "MUJER","QUESO","FUEGO","CASAS","RASGO","GATOS","GOTAS","HUEVO","TINTO","NAVES","VERDE"
$col = 4
foreach ($p in $palabras) {
    $ws.Cells.Item(6, $col).Value = $p
    $col++
}
$ws.Cells.Item(6, 15).Value = "//valores"

$ws.Cells.Item(8, 2).Value = "Información de la estructura:"
$ws.Cells.Item(9, 3).Value = "Tipo: Indexado"
$ws.Cells.Item(10, 3).Value = "Tipos de datos: Almacena valores String"
$ws.Cells.Item(11, 3).Value = "¿Para qué se utilizada?: guarda las palabras que se pueden utilizar para jugar wordix"

# ---------------------------------------------------------------------
# $coleccionPartidas section (new)
# ---------------------------------------------------------------------
$ws.Cells.Item(13, 1).Value = "Representacion del arreglo `$coleccionPartidas"
$ws.Cells.Item(13, 6).Value = "Claves`n|`nv"

$ws.Cells.Item(15, 1).Value = "`$coleccionPartidas="
$indices2 = 0,1,2,3,4,5
$col = 7
foreach ($i in $indices2) {
    $ws.Cells.Item(15, $col).Value = $i
    $col++
}
$ws.Cells.Item(15, 13).Value = "//indices"

$ws.Cells.Item(16, 2).Value = "Información de la estructura:"
$ws.Cells.Item(16, 6).Value = "palabraWordix"
$ws.Cells.Item(16, 7).Value = "QUESO"
$ws.Cells.Item(16, 8).Value = "CASAS"
$ws.Cells.Item(16, 9).Value = "QUESO"
$ws.Cells.Item(16, 13).Value = "//valores"

$ws.Cells.Item(17, 3).Value = "Tipo: Asociativa"
$ws.Cells.Item(17, 6).Value = "jugador"
$ws.Cells.Item(17, 7).Value = "majo"
$ws.Cells.Item(17, 8).Value = "rudolf"
$ws.Cells.Item(17, 9).Value = "pink2000"
$ws.Cells.Item(17, 13).Value = "//valores"

$ws.Cells.Item(18, 3).Value = "Tipos de datos: Almacena valores String e Integer"
$ws.Cells.Item(18, 6).Value = "intentos"
$ws.Cells.Item(18, 7).Value = 0
$ws.Cells.Item(18, 8).Value = 3
$ws.Cells.Item(18, 9).Value = 6
$ws.Cells.Item(18, 13).Value = "//valores"

$ws.Cells.Item(19, 3).Value = "¿Para qué es utilizada?: Guardar los datos de las partidas jugadas"
$ws.Cells.Item(19, 6).Value = "puntaje"
$ws.Cells.Item(19, 7).Value = 0
$ws.Cells.Item(19, 8).Value = 14
$ws.Cells.Item(19, 9).Value = 10
$ws.Cells.Item(19, 13).Value = "//valores"

Write-Output "content written"
